$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Correct Sharon Wu's e-mail address (row 9, column C) and repoint its hyperlink.
$ws.Range("C9").Value = "wuy324@mcmaster.ca"
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:wuy324@mcmaster.ca") | Out-Null

# 2. Correction of the "Second Reviewer" columns (J, K, L) for rows 3-10.
#    Each column cyclically references the next row (wrapping from row 10 back to row 3):
#      J[r] = B[next row]   (next expert, cycling through all 8 project owners)
#      K[r] = J[next row]
#      L[r] = K[next row]
for ($r = 3; $r -le 10; $r++) {
    $nr = $r + 1
    if ($nr -gt 10) { $nr = 3 }
    $ws.Cells.Item($r, 10).Formula = "=B" + $nr
}
for ($r = 3; $r -le 10; $r++) {
    $nr = $r + 1
    if ($nr -gt 10) { $nr = 3 }
    $ws.Cells.Item($r, 11).Formula = "=J" + $nr
}
for ($r = 3; $r -le 10; $r++) {
    $nr = $r + 1
    if ($nr -gt 10) { $nr = 3 }
    $ws.Cells.Item($r, 12).Formula = "=K" + $nr
}

# 3. Move the active selection from D10 to C10.
$ws.Range("C10").Select() | Out-Null
